$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1157.387157691246, 981.181658104533, 605.862092896851, 1292.426469196052, 1165.712211379992, 1157.793623994623, 1295.177072832771),
    @(1208.207581561454, 1001.504143663945, 576.5620703043116, 1316.046994716525, 1212.959743826953, 1167.018325727053, 1317.500688659987),
    @(1063.700906389054, 1039.496240544647, 347.2708387582935, 1258.863885985375, 1073.324544809446, 1102.723238839376, 1261.27490463378),
    @(1281.062412258313, 1131.549665732404, 219.0103028675416, 1356.706564784675, 1282.391924194888, 1176.781486762983, 1356.938493712931),
    @(1262.312376596757, 1122.600416929733, 99.27596812055262, 1348.444573026424, 1264.196820163211, 1140.941785076102, 1348.751454511466),
    @(1305.139994837836, 1176.796944956915, 187.2747926392082, 1364.568010954196, 1306.110061062109, 1207.976052222934, 1364.810771981189),
    @(1187.049016695813, 722.3126791818364, 199.3182319901114, 1250.574057006675, 1187.12235824837, 807.3842191477663, 1250.645859292783),
    @(1317.148662192248, 1065.463844133058, 288.5115740790028, 1354.716259628011, 1318.465130650536, 1118.584399760819, 1355.153956130702),
    @(1315.913788041324, 1157.648520803099, 467.6391110285813, 1364.70517786298, 1318.515521661503, 1256.002499623704, 1366.076505990439),
    @(1302.951897844415, 1160.705684869071, 403.3364748608295, 1365.677207529634, 1305.448640131307, 1248.948808275188, 1366.774754614911),
    @(1157.343514172884, 1172.403438093525, 113.2232041804238, 1344.200000400687, 1162.266216930732, 1192.728605902402, 1344.763941678255),
    @(1308.382254672927, 1155.606495902057, 362.3155122142903, 1362.354721856804, 1310.593135206027, 1237.533746002944, 1363.549160780126)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt $data[$i].Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}
